$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 930.6070999999999
$ws.Range("J17").Value = 1052.591
$ws.Range("L17").Value = 3157.773
$ws.Range("N17").Value = -3493.773

$ws.Range("H93").Value = 54800
$ws.Range("J93").Value = 54800
$ws.Range("L93").Value = 54800
$ws.Range("N93").Value = -59792

$ws.Range("H129").Value = 2028.7878
$ws.Range("I129").Value = 662.55554
$ws.Range("J129").Value = 2541.125
$ws.Range("K129").Value = 1987.66662
$ws.Range("L129").Value = 7623.375
$ws.Range("M129").Value = 3012.33338
$ws.Range("N129").Value = -17623.375

$ws.Range("H138").Value = 1625.79
$ws.Range("I138").Value = 1021.45654
$ws.Range("J138").Value = 2140.5925
$ws.Range("K138").Value = 3064.36962
$ws.Range("L138").Value = 6421.7775
$ws.Range("M138").Value = 2075.63038
$ws.Range("N138").Value = -16701.7775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2711
$ws.Range("I20").Value = 1566.6666
$ws.Range("K20").Value = 1566.6666
$ws.Range("M20").Value = -1319.6666

$ws.Range("H26").Value = 24571.75
$ws.Range("I26").Value = 11611
$ws.Range("K26").Value = 11611
$ws.Range("M26").Value = -11319

$ws.Range("H29").Value = 1001781.6
$ws.Range("I29").Value = 1667296.6
$ws.Range("J29").Value = 3509
$ws.Range("K29").Value = 1667296.6
$ws.Range("L29").Value = 3509
$ws.Range("M29").Value = -1667007.6
$ws.Range("N29").Value = -4087

$ws.Range("H36").Value = 496.4
$ws.Range("I36").Value = 496.4
$ws.Range("K36").Value = 496.4
$ws.Range("M36").Value = 37.60000000000002

$ws.Range("H86").Value = 3425
$ws.Range("I86").Value = 10700
$ws.Range("K86").Value = 10700
$ws.Range("M86").Value = -9577

$ws.Range("H89").Value = 3425
$ws.Range("I89").Value = 10700
$ws.Range("K89").Value = 53500
$ws.Range("M89").Value = -47884

$ws.Range("H109").Value = 30295
$ws.Range("J109").Value = 30295
$ws.Range("L109").Value = 30295
$ws.Range("N109").Value = -33069

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 28571.5
$ws.Range("J28").Value = 28571.5
$ws.Range("L28").Value = 28571.5
$ws.Range("N28").Value = -29061.5

$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250

$ws.Range("H94").Value = 1285.0834
$ws.Range("I94").Value = 930.5
$ws.Range("J94").Value = 1462.375
$ws.Range("K94").Value = 930.5
$ws.Range("L94").Value = 1462.375
$ws.Range("M94").Value = -479.5
$ws.Range("N94").Value = -2364.375

$ws.Range("H98").Value = 63500
$ws.Range("I98").Value = 24000
$ws.Range("J98").Value = 76666.664
$ws.Range("K98").Value = 24000
$ws.Range("L98").Value = 76666.664
$ws.Range("M98").Value = -21754
$ws.Range("N98").Value = -81158.664

$ws.Range("H132").Value = 1728.3529
$ws.Range("I132").Value = 1081
$ws.Range("J132").Value = 2653.1428
$ws.Range("K132").Value = 3243
$ws.Range("L132").Value = 7959.428400000001
$ws.Range("M132").Value = -713
$ws.Range("N132").Value = -13019.4284

$ws.Range("H134").Value = 1463.1666
$ws.Range("I134").Value = 995.76
$ws.Range("K134").Value = 2987.28
$ws.Range("M134").Value = -452.2799999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 816.7586
$ws.Range("I5").Value = 495.1
$ws.Range("J5").Value = 1531.5555
$ws.Range("K5").Value = 1485.3
$ws.Range("L5").Value = 4594.666499999999
$ws.Range("M5").Value = -1373.3
$ws.Range("N5").Value = -4818.666499999999

$ws.Range("H59").Value = 2174.25
$ws.Range("I59").Value = 1200
$ws.Range("J59").Value = 2499
$ws.Range("K59").Value = 3600
$ws.Range("L59").Value = 7497
$ws.Range("M59").Value = -3060
$ws.Range("N59").Value = -8577

$ws.Range("H135").Value = 816.7586
$ws.Range("I135").Value = 495.1
$ws.Range("J135").Value = 1531.5555
$ws.Range("K135").Value = 4455.900000000001
$ws.Range("L135").Value = 13783.9995
$ws.Range("M135").Value = -1920.900000000001
$ws.Range("N135").Value = -18853.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1095
$ws.Range("I3").Value = 2040
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 2040
$ws.Range("L3").Value = 150
$ws.Range("M3").Value = -1924
$ws.Range("N3").Value = -382

$ws.Range("H70").Value = 7310.826
$ws.Range("I70").Value = 9200
$ws.Range("J70").Value = 6303.2666
$ws.Range("K70").Value = 9200
$ws.Range("L70").Value = 6303.2666
$ws.Range("M70").Value = -8930
$ws.Range("N70").Value = -6843.2666

$ws.Range("H73").Value = 7310.826
$ws.Range("I73").Value = 9200
$ws.Range("J73").Value = 6303.2666
$ws.Range("K73").Value = 9200
$ws.Range("L73").Value = 6303.2666
$ws.Range("M73").Value = -8264
$ws.Range("N73").Value = -8175.2666

$ws.Range("H107").Value = 1143.7778
$ws.Range("I107").Value = 1315.6666
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1315.6666
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 604.3334
$ws.Range("N107").Value = -4640

$ws.Range("H132").Value = 2061.1904
$ws.Range("I132").Value = 1194.65
$ws.Range("J132").Value = 3568.2173
$ws.Range("K132").Value = 3583.95
$ws.Range("L132").Value = 10704.6519
$ws.Range("M132").Value = -1053.95
$ws.Range("N132").Value = -15764.6519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 667189.3
$ws.Range("I46").Value = 478
$ws.Range("J46").Value = 1000545
$ws.Range("K46").Value = 478
$ws.Range("L46").Value = 1000545
$ws.Range("M46").Value = -290
$ws.Range("N46").Value = -1000921

$ws.Range("H105").Value = 28000
$ws.Range("J105").Value = 28000
$ws.Range("L105").Value = 28000
$ws.Range("N105").Value = -34988

$ws.Range("H132").Value = 21512.86
$ws.Range("I132").Value = 30326.895
$ws.Range("J132").Value = 3884.7896
$ws.Range("K132").Value = 90980.685
$ws.Range("L132").Value = 11654.3688
$ws.Range("M132").Value = -88450.685
$ws.Range("N132").Value = -16714.3688

$ws.Range("H136").Value = 1718.2667
$ws.Range("I136").Value = 1397.4839
$ws.Range("J136").Value = 2428.5715
$ws.Range("K136").Value = 4192.4517
$ws.Range("L136").Value = 7285.7145
$ws.Range("M136").Value = -1642.4517
$ws.Range("N136").Value = -12385.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1703.9385
$ws.Range("I132").Value = 1440.081
$ws.Range("J132").Value = 2052.6072
$ws.Range("K132").Value = 4320.242999999999
$ws.Range("L132").Value = 6157.821599999999
$ws.Range("M132").Value = -1790.242999999999
$ws.Range("N132").Value = -11217.8216

$ws.Range("H136").Value = 1470.0566
$ws.Range("I136").Value = 1401.027
$ws.Range("J136").Value = 1629.6875
$ws.Range("K136").Value = 4203.081
$ws.Range("L136").Value = 4889.0625
$ws.Range("M136").Value = -1653.081
$ws.Range("N136").Value = -9989.0625
